# SmartStore Working.xlsx update
# Adds Space/Purchase/Sale/Stock-info columns to sheet "02Dec2022",
# removes the old "Spit/StandUp" demo walkthrough content from that sheet,
# adds a new Cycle-counting / Location block, and makes "02Dec2022" the
# active sheet/tab again (it was "Sheet3" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02Dec2022")

# --- New "Input to the System" column (G3:G6) ------------------------------
$ws.Range("G3").Value = "Purchase Information"
$ws.Range("G4").Value = "Sale Informatin"
$ws.Range("G5").Value = "Stock Information"
$ws.Range("G6").Value = "Breakage, Expiry Information"

# G7 used to hold "Margin in Currency not in percentage" - removed entirely
$ws.Range("G7").Clear()

# --- Row 16-18: drop the old F-column demo text, add new G-column values ---
$ws.Range("F16").Clear()
$ws.Range("G16").Value = "Space Information"

$ws.Range("F17").Clear()
$ws.Range("G17").Value = "Shelf Life"

$ws.Range("F18").Clear()
$ws.Range("G18").Value = "Margin in Currency not in percentage"

# --- Row 19: remove old F/G/I/J demo content, add K19 ----------------------
$ws.Range("F19").Clear()
$ws.Range("G19").Clear()
$ws.Range("I19").Clear()
$ws.Range("J19").ClearContents()   # keep the bold style, just blank the text
$ws.Range("K19").Value = "Sale Information"

# --- Row 20: remove old F/G/I/J demo content, add K20 ----------------------
$ws.Range("F20").Clear()
$ws.Range("G20").ClearContents()   # keep the bold style, just blank the text
$ws.Range("I20").Clear()
$ws.Range("J20").ClearContents()   # keep the bold style, just blank the text
$ws.Range("K20").Value = "Bar / QR Code"

# --- Row 21: remove old F/I/J demo content, add K21 -------------------------
$ws.Range("F21").Clear()
$ws.Range("I21").Clear()
$ws.Range("J21").Clear()
$ws.Range("K21").Value = "Stock Level"

# --- Row 22: remove old I/J demo content, add K22 ---------------------------
$ws.Range("I22").Clear()
$ws.Range("J22").Clear()
$ws.Range("K22").Value = "Wastage, breakage. Expiry"

# --- Row 23: remove old I/J demo content (row becomes empty) ---------------
$ws.Range("I23").Clear()
$ws.Range("J23").Clear()

# --- Row 24: remove old I/J demo content, add K24 ---------------------------
$ws.Range("I24").Clear()
$ws.Range("J24").Clear()
$ws.Range("K24").Value = "Stock Information at the time of Purchase Entry"

# --- Row 25: remove old I/J demo content, add C25/K25 -----------------------
$ws.Range("I25").Clear()
$ws.Range("J25").Clear()
$ws.Range("C25").Value = "I have huge number of items, it is very difficult for me to enter and maintain these in the system"
$ws.Range("K25").Value = "Cycle counting"

# --- New rows 26-29: Location list in column L ------------------------------
$ws.Range("L26").Value = "Location Manual"
$ws.Range("L27").Value = "Random"
$ws.Range("L28").Value = "Location Sequence"
$ws.Range("L29").Value = "Location Random"

# --- Make "02Dec2022" the active sheet/tab again, with C24 selected --------
$ws.Activate()
$ws.Range("C24").Select()
